$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simple Data")
$ws.Activate()

# Labels for the per-column averages (columns D..H -> rows 2..6 in col J)
$ws.Range("J3").Value = "avg2"
$ws.Range("J2").Value = "avg1"
$ws.Range("J4").Value = "avg3"
$ws.Range("J5").Value = "avg4"
$ws.Range("J6").Value = "avg5"

# Average formulas for columns D..H, computed over rows 2:68
$ws.Range("K2").Formula = "=AVERAGE(D2:D68)"
$ws.Range("K3").Formula = "=AVERAGE(E2:E68)"
$ws.Range("K4").Formula = "=AVERAGE(F2:F68)"
$ws.Range("K5").Formula = "=AVERAGE(G2:G68)"
$ws.Range("K6").Formula = "=AVERAGE(H2:H68)"

# Update the active selection to match the recorded UI state after the edit
$ws.Range("K7").Select()
